# Update absenteeism data rows 2-11 with new values, preserving formatting/styles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 11789; B = "Marina Costa";          C = "Vendas";                  D = "Doença";              E = 2; F = 45081; G = 8178.01 },
    @{ Row = 3;  A = 9953;  B = "Rebeca Martins";         C = "P&D";                     D = "Viagem de negócios";  E = 2; F = 45094; G = 12439.21 },
    @{ Row = 4;  A = 78014; B = "Lucas Gabriel da Mota";  C = "TI";                      D = "Viagem de negócios";  E = 7; F = 45079; G = 10532.14 },
    @{ Row = 5;  A = 33330; B = "Luiz Miguel Pires";      C = "Financeiro";              D = "Viagem de negócios";  E = 1; F = 45089; G = 7823.31 },
    @{ Row = 6;  A = 65128; B = "Ana Beatriz Dias";       C = "Operações";               D = "Problemas pessoais";  E = 6; F = 45097; G = 6621.36 },
    @{ Row = 7;  A = 25489; B = "Nicole Caldeira";        C = "Recursos Humanos";        D = "Problemas pessoais";  E = 7; F = 45089; G = 5861.1 },
    @{ Row = 8;  A = 17009; B = "Heitor Vieira";          C = "Recursos Humanos";        D = "Viagem de negócios";  E = 3; F = 45087; G = 5645.18 },
    @{ Row = 9;  A = 67378; B = "Pietro Alves";           C = "Marketing";               D = "Consulta médica";     E = 1; F = 45100; G = 3782.74 },
    @{ Row = 10; A = 1466;  B = "Maitê Silva";            C = "Atendimento ao Cliente";  D = "Consulta médica";     E = 5; F = 45105; G = 7031.39 },
    @{ Row = 11; A = 8824;  B = "Fernanda Pinto";         C = "Recursos Humanos";        D = "Viagem de negócios";  E = 1; F = 45104; G = 11673.54 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
